# ADD: stl file for stage body
#
# Row 11 ("Body") is split into two rows: "Body_base" (existing row 11) and
# a new "Body_table" row (new row 12). Both reference a "3D printer "
# (trailing space) manufacturing note, matching the new row's plain
# "3D printer" text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing "Body" row item to "Body_base".
$ws.Range("B11").Value = "Body_base"

# Insert the new row for the table/lid part the STL addition is for.
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Body_table"
$ws.Range("C12").Value = "-"
$ws.Range("D12").Value = "3D printer"

# The existing row's manufacturing note gets a trailing space to
# distinguish it from the new row's "3D printer" note.
$ws.Range("D11").Value = "3D printer "

# Row 7 no longer needs an explicit custom height - let Excel use the
# default row height again.
$ws.Rows.Item(7).AutoFit()

# Update the active selection to mirror where the author left off editing.
$ws.Range("D11").Select()
